# ZBP_08_pocet_aktivit.xlsx update: add a new data-collection wave
# (22. 2. 2022) as a new trailing column on both sheets, and correct a
# handful of previously-provisional values in the now-second-to-last
# column now that final data is in.
#
# Sheet "data": new column AO (after AN), rows 1-76.
# Sheet "pocetR": new column AN (after AM), rows 1-26.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "data" - add column AO
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# Header cell AO1: copy the date-header formatting from AN1, then set text
$srcHdr1 = $ws1.Range("AN1")
$dstHdr1 = $ws1.Range("AO1")
$srcHdr1.Copy()
$dstHdr1.PasteSpecial(-4122)
$excel.CutCopyMode = 0
$dstHdr1.Value = "22. 2. 2022"

# New AO values for rows 2-76 (no special formatting on data cells)
$ao1Data = @{
    2 = 0.64
    3 = 0.23
    4 = 0.13
    5 = 0.72
    6 = 0.14
    7 = 0.14
    8 = 0.61
    9 = 0.25
    10 = 0.14
    11 = 0.78
    12 = 0.13
    13 = 0.09
    14 = 0.73
    15 = 0.18
    16 = 0.09
    17 = 0.5600000000000001
    18 = 0.16
    19 = 0.28
    20 = 0.47
    21 = 0.34
    22 = 0.19
    23 = 0.8100000000000001
    24 = 0.12
    25 = 0.07000000000000001
    26 = 0.6899999999999999
    27 = 0.21
    28 = 0.1
    29 = 0.47
    30 = 0.32
    31 = 0.21
    32 = 0.42
    33 = 0.32
    34 = 0.26
    35 = 0.6
    36 = 0.24
    37 = 0.16
    38 = 0.73
    39 = 0.2
    40 = 0.07000000000000001
    41 = 0.79
    42 = 0.15
    43 = 0.06
    44 = 0.64
    45 = 0.22
    46 = 0.14
    47 = 0.43
    48 = 0.34
    49 = 0.23
    50 = 0.71
    51 = 0.19
    52 = 0.1
    53 = 0.54
    54 = 0.31
    55 = 0.15
    56 = 0.63
    57 = 0.22
    58 = 0.15
    59 = 0.65
    60 = 0.23
    61 = 0.12
    62 = 0.62
    63 = 0.23
    64 = 0.15
    65 = 0.65
    66 = 0.22
    67 = 0.13
    68 = 0.62
    69 = 0.24
    70 = 0.14
    71 = 0.62
    72 = 0.22
    73 = 0.16
    74 = 0.6
    75 = 0.25
    76 = 0.15
}
foreach ($r in $ao1Data.Keys) {
    $ws1.Range("AO" + $r).Value = $ao1Data[$r]
}

# Revised AN values (previously-provisional figures corrected now that
# the 22. 2. 2022 wave has landed)
$an1Updates = @{
    3 = 0.25
    4 = 0.17
    6 = 0.15
    7 = 0.14
    8 = 0.54
    10 = 0.18
    15 = 0.23
    16 = 0.11
    17 = 0.58
    18 = 0.19
    19 = 0.23
    20 = 0.39
    22 = 0.25
    23 = 0.77
    25 = 0.1
    26 = 0.63
    27 = 0.25
    35 = 0.52
    37 = 0.17
    41 = 0.71
    43 = 0.11
    44 = 0.62
    46 = 0.15
    47 = 0.42
    49 = 0.23
    50 = 0.62
    52 = 0.14
    59 = 0.61
    60 = 0.24
    63 = 0.27
    64 = 0.18
    65 = 0.61
    67 = 0.17
    68 = 0.55
    70 = 0.16
    72 = 0.31
    73 = 0.15
}
foreach ($r in $an1Updates.Keys) {
    $ws1.Range("AN" + $r).Value = $an1Updates[$r]
}

# Footer title row: bump the "aktualizace" (update) date
$ws1.Range("A77").Value = "Život během pandemie, Počet protektivních aktivit, % respondentů celkově a ve skupinách, aktualizace 2. 3. 2022"

# ---------------------------------------------------------------------
# Sheet 2: "pocetR" - add column AN
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

# Header cell AN1: copy the date-header formatting from AM1, then set text
$srcHdr2 = $ws2.Range("AM1")
$dstHdr2 = $ws2.Range("AN1")
$srcHdr2.Copy()
$dstHdr2.PasteSpecial(-4122)
$excel.CutCopyMode = 0
$dstHdr2.Value = "22. 2. 2022"

# New AN values for rows 2-26 (no special formatting on data cells)
$an2Data = @{
    2 = 1786
    3 = 366
    4 = 1420
    5 = 275
    6 = 808
    7 = 91
    8 = 612
    9 = 445
    10 = 663
    11 = 678
    12 = 294
    13 = 528
    14 = 703
    15 = 567
    16 = 778
    17 = 441
    18 = 682
    19 = 421
    20 = 683
    21 = 863
    22 = 923
    23 = 931
    24 = 413
    25 = 206
    26 = 236
}
foreach ($r in $an2Data.Keys) {
    $ws2.Range("AN" + $r).Value = $an2Data[$r]
}

# Revised AM values (previously-provisional figures corrected now that
# the 22. 2. 2022 wave has landed)
$am2Updates = @{
    2 = 1848
    3 = 419
    4 = 1429
    5 = 305
    6 = 815
    7 = 114
    8 = 614
    9 = 451
    10 = 689
    11 = 708
    13 = 546
    14 = 734
    15 = 499
    16 = 773
    17 = 576
    18 = 683
    19 = 504
    20 = 661
    21 = 893
    22 = 955
    23 = 973
    24 = 416
    25 = 215
    26 = 244
}
foreach ($r in $am2Updates.Keys) {
    $ws2.Range("AM" + $r).Value = $am2Updates[$r]
}

# Footer title row: bump the "aktualizace" (update) date
$ws2.Range("A27").Value = "Život během pandemie, Počet protektivních aktivit, velikost dotázaného souboru celkově a ve skupinách, aktualizace 2. 3. 2022"
